# Week 34 questions and profiles updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 85: answer flipped from 1 to 0 -------------------------------
$ws.Range("C85").Value = 0

# --- New profile rows 105-108 ------------------------------------------
# Row 105
$ws.Range("C97").Copy()
$ws.Range("C105").PasteSpecial(-4122)
$ws.Range("A104").Copy()
$ws.Range("A105").PasteSpecial(-4122)
$ws.Range("A105").Value = "Harsh Mehta"
$ws.Range("B105").Value = "harsh_mehta2"
$ws.Range("C105").Value = 0

# Row 106
$ws.Range("C90").Copy()
$ws.Range("C106").PasteSpecial(-4122)
$ws.Range("A104").Copy()
$ws.Range("A106").PasteSpecial(-4122)
$ws.Range("A106").Value = "Priya Agrawal"
$ws.Range("B106").Value = "priya_agrawal3"
$ws.Range("C106").Value = 0

# Row 107
$ws.Range("C97").Copy()
$ws.Range("C107").PasteSpecial(-4122)
$ws.Range("A104").Copy()
$ws.Range("A107").PasteSpecial(-4122)
$ws.Range("A107").Value = "Ravi Kumar Sharma"
$ws.Range("B107").Value = "rksharma2180"
$ws.Range("B107").WrapText = $false
$ws.Range("C107").Value = 0

# Row 108
$ws.Range("C90").Copy()
$ws.Range("C108").PasteSpecial(-4122)
$ws.Range("A104").Copy()
$ws.Range("A108").PasteSpecial(-4122)
$ws.Range("A108").Value = "Nimit Bansal"
$ws.Range("B108").Value = "NimitBnsl"
$ws.Range("B108").WrapText = $false
$ws.Range("C108").Value = 0

# --- View state: scroll + selection on B108 -----------------------------
$ws.Range("B108").Select()
$excel.ActiveWindow.ScrollRow = 74
$excel.ActiveWindow.ScrollColumn = 1
